$wb = $excel.ActiveWorkbook

# --- "WithHeadings" sheet: shrink the selection from A3:XFD4 down to just A3 ---
$wsHeadings = $wb.Worksheets.Item("WithHeadings")
$wsHeadings.Range("A3").Select() | Out-Null

# --- "HeadingsOnRowThree" sheet: update A4, drop the stray K21 "six" entry ---
$wsHeadingsRow3 = $wb.Worksheets.Item("HeadingsOnRowThree")
$wsHeadingsRow3.Range("A4").Value = 1234
$wsHeadingsRow3.Range("K21").ClearContents()

# Make this the active sheet / selection (D13), which also clears tabSelected
# on whichever sheet previously held it ("TypeTests").
$wsHeadingsRow3.Activate()
$wsHeadingsRow3.Range("D13").Select() | Out-Null
